$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Insl5"
$ws.Range("C2").Value = "Rxfp4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1188713333333333
$ws.Range("H2").Value = 0.356614
$ws.Range("I2").Value = 0.2549271348773238
$ws.Range("J2").Value = 0.339160140832479
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.442313
$ws.Range("N2").Value = 10.326939
$ws.Range("O2").Value = 0.2120737065114005
$ws.Range("P2").Value = 0.2368526181325179
$ws.Range("Q2").Value = 0.4091923360606666
$ws.Range("R2").Value = 3.682731024546
$ws.Range("S2").Value = 0.05406334238376578
$ws.Range("T2").Value = 0.08033096732236612

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Insl5"
$ws.Range("C3").Value = "Rxfp4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1188713333333333
$ws.Range("H3").Value = 0.356614
$ws.Range("I3").Value = 0.2549271348773238
$ws.Range("J3").Value = 0.339160140832479
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.571969666666666
$ws.Range("N3").Value = 7.715909
$ws.Range("O3").Value = 0.158453673516874
$ws.Range("P3").Value = 0.1769675649214407
$ws.Range("Q3").Value = 0.3057334635695555
$ws.Range("R3").Value = 2.751601172126
$ws.Range("S3").Value = 0.04039414100044358
$ws.Range("T3").Value = 0.06002034424153669

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Insl5"
$ws.Range("C4").Value = "Rxfp4"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1188713333333333
$ws.Range("H4").Value = 0.356614
$ws.Range("I4").Value = 0.2549271348773238
$ws.Range("J4").Value = 0.339160140832479
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.535264
$ws.Range("N4").Value = 7.605791999999999
$ws.Range("O4").Value = 0.1561923141402073
$ws.Range("P4").Value = 0.174441985971967
$ws.Range("Q4").Value = 0.301370212032
$ws.Range("R4").Value = 2.712331908288
$ws.Range("S4").Value = 0.03981765913362195
$ws.Range("T4").Value = 0.05916376852934966

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Insl5"
$ws.Range("C5").Value = "Rxfp4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1188713333333333
$ws.Range("H5").Value = 0.356614
$ws.Range("I5").Value = 0.2549271348773238
$ws.Range("J5").Value = 0.339160140832479
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.0943505
$ws.Range("N5").Value = 10.188701
$ws.Range("O5").Value = 0.3138522826957358
$ws.Range("P5").Value = 0.2336820724146239
$ws.Range("Q5").Value = 0.6055722364023333
$ws.Range("R5").Value = 3.633433418414
$ws.Range("S5").Value = 0.08000946320233182
$ws.Range("T5").Value = 0.07925564459016937

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Insl5"
$ws.Range("C6").Value = "Rxfp4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1188713333333333
$ws.Range("H6").Value = 0.356614
$ws.Range("I6").Value = 0.2549271348773238
$ws.Range("J6").Value = 0.339160140832479
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.587785
$ws.Range("N6").Value = 7.763355
$ws.Range("O6").Value = 0.1594280231357824
$ws.Range("P6").Value = 0.1780557585594505
$ws.Range("Q6").Value = 0.30761345333
$ws.Range("R6").Value = 2.76852107997
$ws.Range("S6").Value = 0.0406425291571607
$ws.Range("T6").Value = 0.06038941614905712

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Insl5"
$ws.Range("C7").Value = "Rxfp4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.347424
$ws.Range("H7").Value = 0.694848
$ws.Range("I7").Value = 0.7450728651226762
$ws.Range("J7").Value = 0.6608398591675211
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.442313
$ws.Range("N7").Value = 10.326939
$ws.Range("O7").Value = 0.2120737065114005
$ws.Range("P7").Value = 0.2368526181325179
$ws.Range("Q7").Value = 1.195942151712
$ws.Range("R7").Value = 7.175652910271999
$ws.Range("S7").Value = 0.1580103641276347
$ws.Range("T7").Value = 0.1565216508101518

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Insl5"
$ws.Range("C8").Value = "Rxfp4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.347424
$ws.Range("H8").Value = 0.694848
$ws.Range("I8").Value = 0.7450728651226762
$ws.Range("J8").Value = 0.6608398591675211
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.571969666666666
$ws.Range("N8").Value = 7.715909
$ws.Range("O8").Value = 0.158453673516874
$ws.Range("P8").Value = 0.1769675649214407
$ws.Range("Q8").Value = 0.893563989472
$ws.Range("R8").Value = 5.361383936832
$ws.Range("S8").Value = 0.1180595325164305
$ws.Range("T8").Value = 0.116947220679904

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Insl5"
$ws.Range("C9").Value = "Rxfp4"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.347424
$ws.Range("H9").Value = 0.694848
$ws.Range("I9").Value = 0.7450728651226762
$ws.Range("J9").Value = 0.6608398591675211
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.535264
$ws.Range("N9").Value = 7.605791999999999
$ws.Range("O9").Value = 0.1561923141402073
$ws.Range("P9").Value = 0.174441985971967
$ws.Range("Q9").Value = 0.8808115599359999
$ws.Range("R9").Value = 5.284869359616
$ws.Range("S9").Value = 0.1163746550065853
$ws.Range("T9").Value = 0.1152782174426174

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Insl5"
$ws.Range("C10").Value = "Rxfp4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.347424
$ws.Range("H10").Value = 0.694848
$ws.Range("I10").Value = 0.7450728651226762
$ws.Range("J10").Value = 0.6608398591675211
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.0943505
$ws.Range("N10").Value = 10.188701
$ws.Range("O10").Value = 0.3138522826957358
$ws.Range("P10").Value = 0.2336820724146239
$ws.Range("Q10").Value = 1.769899628112
$ws.Range("R10").Value = 7.079598512448
$ws.Range("S10").Value = 0.233842819493404
$ws.Range("T10").Value = 0.1544264278244545

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Insl5"
$ws.Range("C11").Value = "Rxfp4"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.347424
$ws.Range("H11").Value = 0.694848
$ws.Range("I11").Value = 0.7450728651226762
$ws.Range("J11").Value = 0.6608398591675211
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.587785
$ws.Range("N11").Value = 7.763355
$ws.Range("O11").Value = 0.1594280231357824
$ws.Range("P11").Value = 0.1780557585594505
$ws.Range("Q11").Value = 0.8990586158399999
$ws.Range("R11").Value = 5.39435169504
$ws.Range("S11").Value = 0.1187854939786217
$ws.Range("T11").Value = 0.1176663424103934
